$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume/Number and report week dates) ---
$ws.Range("A8").Value = "Volume 32   Number  11"
$ws.Range("C9").Value = "Report Covering the Week  3/10/2025  Through  3/16/2025"

# --- Cells that flip between numeric and text representation ---
# Step 1: force the text-bearing cells to store literal text (shared string)
# by setting an explicit Text number format before assigning the value.
$ws.Range("F27").NumberFormat = "@"
$ws.Range("F27").Value = "0"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "***.*"

# Step 2: copy the correct number format from a same-styled neighbour cell
# onto each changed cell (format only -- does not disturb the stored value/type).
$ws.Range("C16").Copy()
$ws.Range("D15").PasteSpecial(-4122)
$ws.Range("D33").PasteSpecial(-4122)
$ws.Range("G33").PasteSpecial(-4122)
$ws.Range("H15").Copy()
$ws.Range("E15").PasteSpecial(-4122)
$ws.Range("E33").PasteSpecial(-4122)
$ws.Range("H33").PasteSpecial(-4122)
$ws.Range("C14").Copy()
$ws.Range("F27").PasteSpecial(-4122)
$ws.Range("D31").PasteSpecial(-4122)
$ws.Range("E14").Copy()
$ws.Range("E31").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Step 3: now assign the final numeric values to the cells that end up numeric.
$ws.Range("D15").Value = 2
$ws.Range("E15").Value = -100
$ws.Range("D33").Value = 1
$ws.Range("E33").Value = -100
$ws.Range("G33").Value = 1
$ws.Range("H33").Value = 0

# --- Plain value updates (style/type unchanged) ---
$ws.Range("N14").Value = -60
$ws.Range("G15").Value = 3
$ws.Range("J15").Value = 6
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 20
$ws.Range("C16").Value = 6
$ws.Range("D16").Value = 6
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 23
$ws.Range("G16").Value = 27
$ws.Range("H16").Value = -14.814814814814
$ws.Range("I16").Value = 48
$ws.Range("J16").Value = 72
$ws.Range("K16").Value = -33.333333333333
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -25
$ws.Range("N16").Value = -81.609195402298
$ws.Range("C17").Value = 10
$ws.Range("D17").Value = 16
$ws.Range("E17").Value = -37.5
$ws.Range("F17").Value = 36
$ws.Range("G17").Value = 42
$ws.Range("H17").Value = -14.285714285714
$ws.Range("I17").Value = 96
$ws.Range("J17").Value = 103
$ws.Range("K17").Value = -6.796116504854
$ws.Range("L17").Value = 1.052631578947
$ws.Range("M17").Value = 74.545454545454
$ws.Range("N17").Value = 43.283582089552
$ws.Range("C18").Value = 1
$ws.Range("D18").Value = 1
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 9
$ws.Range("G18").Value = 16
$ws.Range("H18").Value = -43.75
$ws.Range("I18").Value = 31
$ws.Range("J18").Value = 45
$ws.Range("K18").Value = -31.111111111111
$ws.Range("L18").Value = 29.166666666666
$ws.Range("M18").Value = -49.180327868852
$ws.Range("N18").Value = -94.139886578449
$ws.Range("C19").Value = 8
$ws.Range("D19").Value = 16
$ws.Range("E19").Value = -50
$ws.Range("F19").Value = 50
$ws.Range("G19").Value = 60
$ws.Range("H19").Value = -16.666666666666
$ws.Range("I19").Value = 140
$ws.Range("J19").Value = 200
$ws.Range("K19").Value = -30
$ws.Range("L19").Value = -14.634146341463
$ws.Range("M19").Value = 55.555555555555
$ws.Range("N19").Value = -44.444444444444
$ws.Range("C20").Value = 7
$ws.Range("D20").Value = 4
$ws.Range("E20").Value = 75
$ws.Range("F20").Value = 18
$ws.Range("G20").Value = 15
$ws.Range("H20").Value = 20
$ws.Range("I20").Value = 38
$ws.Range("J20").Value = 48
$ws.Range("K20").Value = -20.833333333333
$ws.Range("L20").Value = -40.625
$ws.Range("M20").Value = -17.391304347826
$ws.Range("N20").Value = -91.441441441441
$ws.Range("C21").Value = 32
$ws.Range("D21").Value = 45
$ws.Range("E21").Value = -28.888888888888
$ws.Range("F21").Value = 136
$ws.Range("G21").Value = 163
$ws.Range("H21").Value = -16.564417177914
$ws.Range("I21").Value = 361
$ws.Range("J21").Value = 474
$ws.Range("K21").Value = -23.839662447257
$ws.Range("L21").Value = -9.75
$ws.Range("M21").Value = 13.166144200627
$ws.Range("N21").Value = -76.947637292464
$ws.Range("C22").Value = 3
$ws.Range("E22").Value = 0
$ws.Range("F22").Value = 5
$ws.Range("G22").Value = 7
$ws.Range("H22").Value = -28.571428571428
$ws.Range("I22").Value = 13
$ws.Range("J22").Value = 16
$ws.Range("K22").Value = -18.75
$ws.Range("L22").Value = -27.777777777777
$ws.Range("M22").Value = 85.714285714285
$ws.Range("C24").Value = 28
$ws.Range("D24").Value = 37
$ws.Range("E24").Value = -24.324324324324
$ws.Range("F24").Value = 103
$ws.Range("G24").Value = 197
$ws.Range("H24").Value = -47.715736040609
$ws.Range("I24").Value = 282
$ws.Range("J24").Value = 513
$ws.Range("K24").Value = -45.029239766081
$ws.Range("L24").Value = -40.127388535031
$ws.Range("M24").Value = 42.424242424242
$ws.Range("C25").Value = 11
$ws.Range("D25").Value = 18
$ws.Range("E25").Value = -38.888888888888
$ws.Range("F25").Value = 43
$ws.Range("G25").Value = 113
$ws.Range("H25").Value = -61.946902654867
$ws.Range("I25").Value = 128
$ws.Range("J25").Value = 307
$ws.Range("K25").Value = -58.306188925081
$ws.Range("L25").Value = -50
$ws.Range("C26").Value = 30
$ws.Range("D26").Value = 20
$ws.Range("E26").Value = 50
$ws.Range("F26").Value = 88
$ws.Range("G26").Value = 95
$ws.Range("H26").Value = -7.368421052631
$ws.Range("I26").Value = 223
$ws.Range("J26").Value = 226
$ws.Range("K26").Value = -1.327433628318
$ws.Range("L26").Value = 30.409356725146
$ws.Range("M26").Value = 18.617021276595
$ws.Range("D27").Value = 2
$ws.Range("G27").Value = 4
$ws.Range("H27").Value = -100
$ws.Range("J27").Value = 13
$ws.Range("K27").Value = -38.461538461538
$ws.Range("L27").Value = -20
$ws.Range("C28").Value = 3
$ws.Range("D28").Value = 4
$ws.Range("E28").Value = -25
$ws.Range("G28").Value = 11
$ws.Range("H28").Value = -9.090909090909
$ws.Range("I28").Value = 20
$ws.Range("J28").Value = 29
$ws.Range("K28").Value = -31.034482758620
$ws.Range("J33").Value = 3
$ws.Range("K33").Value = -66.666666666666
